$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do List")
$ws.Range("G11").Formula = "=TODAY()-1"
$ws.Range("D12").Formula = "=TODAY()-1"
$ws.Range("G12").Formula = "=TODAY()-1"
$ws.Range("D13").Formula = "=TODAY()-1"
$ws.Range("G13").Formula = "=TODAY()-1"
$ws.Range("D15").Formula = "=TODAY()-1"
$ws.Range("G15").Formula = "=TODAY()-1"
$ws.Range("B7").Value = 1
$ws.Range("B17").Value = 1
